$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.031.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -7.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.409.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -8.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9975"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9969"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "273.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.66%  "
$ws.Range("E7").Value = "  -6.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3129"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.73"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.010"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06513"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.491"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.165"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -7.35%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.408.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.59%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001015"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.05694"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -13.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -15.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9969"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.570"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -9.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.76"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.266"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "19.970.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.241"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "136.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -10.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.571.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -8.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "109.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.107"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -15.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.317"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -13.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8220"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -15.93%  "
$ws.Range("E34").Value = "  -5.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.397"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.472"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05790"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.820"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9977"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02074"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -7.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1903"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.59%  "
$ws.Range("E43").Value = "  -8.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5270"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -10.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.503"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5134"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.51%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "111.86"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.766"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.035"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -11.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.28%  "

Write-Output "edit complete"
